$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update sheet1 (AddCustomerTest) header/alerttext/amount cells first ---
$ws1.Range("D1").Value = "alerttext"
$ws1.Range("C2").Value = 411033
$ws1.Range("D2").Value = "Customer added successfully"
$ws1.Range("D3").Value = "Customer added successfully"

# --- Add new sheet (OpenAccountTest) after sheet1, populate it next so its
# shared strings land before the row-3 customer name strings below ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "OpenAccountTest"
$ws2.Range("A1").Value = "customer"
$ws2.Range("B1").Value = "currency"
$ws2.Range("A2").Value = "sandip thopate"
$ws2.Range("B2").Value = "rupee"
$ws2.Columns.Item(1).ColumnWidth = 15.5
$ws2.Columns.Item(2).ColumnWidth = 12.6666666667
[void]($ws2.Range("H17").Select())

# --- finish sheet1 row 3 (new shared strings appended last) ---
$ws1.Range("A3").Value = "vishal"
$ws1.Range("B3").Value = "sahu"
$ws1.Range("C3").Value = 411034

# OpenAccountTest becomes the active/selected tab, as in the target workbook
[void]($ws2.Activate())
